# Fruta / hortaliza, semanal
# Add a new week of price data (2023-07-25) for "Feria Lagunitas de Puerto Montt - Kiwi"
# at the top of the block of rows that starts at row 558, pushing the existing
# rows 558..617 down by two rows (to 560..619) and writing the two brand-new
# rows (for quality grades "Primera" and "Segunda") into rows 558 and 559.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 558
$lastRow  = 617
$numNewRows = 2
$lastCol = 20   # column T

# 1) Snapshot the existing block (rows 558..617, columns A..T) before we
#    overwrite anything, so the shift doesn't clobber data we still need.
$data = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    for ($c = 1; $c -le $lastCol; $c++) {
        $data["$r,$c"] = $ws.Cells.Item($r, $c).Value2
    }
}

# 2) Shift the snapshot down by $numNewRows: the row that used to be at
#    position r now lives at r + numNewRows. Write from the bottom up so we
#    never need to re-read a cell we've already overwritten (not strictly
#    required since we work off the snapshot, but keeps intent clear).
for ($r = ($lastRow + $numNewRows); $r -ge ($firstRow + $numNewRows); $r--) {
    $src = $r - $numNewRows
    for ($c = 1; $c -le $lastCol; $c++) {
        $ws.Cells.Item($r, $c).Value = $data["$src,$c"]
    }
}

# 3) Write the two brand-new rows for the new reporting week (2023-07-25,
#    serial 45132) into the now-vacated rows 558 and 559.
$ws.Cells.Item(558, 1).Value  = 4
$ws.Cells.Item(558, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(558, 3).Value  = "Los Lagos"
$ws.Cells.Item(558, 4).Value  = 45132
$ws.Cells.Item(558, 5).Value  = 10
$ws.Cells.Item(558, 6).Value  = "Fruta"
$ws.Cells.Item(558, 7).Value  = 100101
$ws.Cells.Item(558, 8).Value  = "Berries"
$ws.Cells.Item(558, 9).Value  = 100101007
$ws.Cells.Item(558, 10).Value = "Kiwi"
$ws.Cells.Item(558, 11).Value = "Hayward"
$ws.Cells.Item(558, 12).Value = "Primera"
$ws.Cells.Item(558, 13).Value = 300
$ws.Cells.Item(558, 14).Value = 18000
$ws.Cells.Item(558, 15).Value = 18000
$ws.Cells.Item(558, 16).Value = 18000
$ws.Cells.Item(558, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(558, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(558, 19).Value = 1200
$ws.Cells.Item(558, 20).Value = 15

$ws.Cells.Item(559, 1).Value  = 4
$ws.Cells.Item(559, 2).Value  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(559, 3).Value  = "Los Lagos"
$ws.Cells.Item(559, 4).Value  = 45132
$ws.Cells.Item(559, 5).Value  = 10
$ws.Cells.Item(559, 6).Value  = "Fruta"
$ws.Cells.Item(559, 7).Value  = 100101
$ws.Cells.Item(559, 8).Value  = "Berries"
$ws.Cells.Item(559, 9).Value  = 100101007
$ws.Cells.Item(559, 10).Value = "Kiwi"
$ws.Cells.Item(559, 11).Value = "Hayward"
$ws.Cells.Item(559, 12).Value = "Segunda"
$ws.Cells.Item(559, 13).Value = 300
$ws.Cells.Item(559, 14).Value = 15000
$ws.Cells.Item(559, 15).Value = 15000
$ws.Cells.Item(559, 16).Value = 15000
$ws.Cells.Item(559, 17).Value = "$/caja 15 kilos"
$ws.Cells.Item(559, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(559, 19).Value = 1000
$ws.Cells.Item(559, 20).Value = 15

# 4) Apply the same date-number-format style ("s=2" in the original sheet)
#    to the whole D column of the block (558..619), matching every other
#    row in D - this also covers rows 618/619 which are brand-new cells
#    that didn't inherit any formatting from the shift-by-value write.
$dateFormat = $ws.Range("D2").NumberFormat
$ws.Range("D558:D619").NumberFormat = $dateFormat
